$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status: "In Translation" -> "Ready for handoff" (shared across Overview!B2/C2, zh-cn!C2, de-de!C2)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Overview Latest Handoff Date D2
$overview.Range("D2").Value = "2016-26-20 18:26:56"

# zh-cn Latest Handoff Datetime E2
$zhcn.Range("E2").Value = "2016-03-20 18:26:53"

# de-de Latest Handoff Datetime E2
$dede.Range("E2").Value = "2016-03-20 18:26:56"
